$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 (pushes existing rows 6-10 down to 7-11)
$ws.Rows.Item(6).Insert()

# Grow the "Responses" table/autofilter to cover the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F11"))

# Restore the calculated "S.N." column formula for the newly inserted row
$ws.Range("A6").Formula = "=ROWS(`$A`$2:Responses[[#This Row],[RESPONSE NAME]])"

# The row insertion can leave the table's last-row calculated formula in a
# "current row" relative form that evaluates incorrectly - rewrite it so it
# stays the normal structured reference bound to the table.
$ws.Range("A11").Formula = "=ROWS(`$A`$2:Responses[[#This Row],[RESPONSE NAME]])"

# Populate the new row: TEXT (D) first, then RESPONSE NAME (C), for the
# fallback utterance - matches the order the shared strings were authored in.
$ws.Range("D6").Value = "I'm sorry, I didn't quite understand that. could you rephrase your query?"
$ws.Range("C6").Value = "utter_please_rephrase"

# Update the view: scroll down a bit and select the newly typed cell
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("D6").Select()
